# Auto-generated edit script for cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. '63.327.18', '5.80') that must stay
# literal text, not be auto-coerced to a number by Excel's input parser. A leading
# apostrophe (quote-prefix) forces text entry exactly like a user typing '5.80 would.
$ws.Range("D2").Value = '''63.327.18'
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("D3").Value = '''3.060.76'
$ws.Range("E3").Value = '  -3.59%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '''587.39'
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").Value = '''154.43'
$ws.Range("E6").Value = '  +4.00%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '''0.535'
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").Value = '''3.059.16'
$ws.Range("E9").Value = '  -3.24%  '
$ws.Range("D10").Value = '''0.155'
$ws.Range("E10").Value = '  -4.54%  '
$ws.Range("D11").Value = '''5.80'
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = '''0.448'
$ws.Range("E12").Value = '  -2.83%  '
$ws.Range("D13").Value = '''36.71'
$ws.Range("E13").Value = '  -2.66%  '
$ws.Range("D14").Value = '''0.0000236'
$ws.Range("E14").Value = '  -4.99%  '
$ws.Range("E15").Value = '  -2.29%  '
$ws.Range("D16").Value = '''3.573.85'
$ws.Range("E16").Value = '  -3.30%  '
$ws.Range("D17").Value = '''63.360.62'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '''7.09'
$ws.Range("E18").Value = '  -3.08%  '
$ws.Range("D19").Value = '''3.064.72'
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").Value = '''468.86'
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").Value = '''14.24'
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("D22").Value = '''0.701'
$ws.Range("E22").Value = '  -5.13%  '
$ws.Range("D23").Value = '''7.46'
$ws.Range("E23").Value = '  -3.21%  '
$ws.Range("D24").Value = '''2.38'
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").Value = '''80.37'
$ws.Range("E25").Value = '  -1.55%  '
$ws.Range("D26").Value = '''12.74'
$ws.Range("E26").Value = '  -4.10%  '
$ws.Range("D27").Value = '''10.28'
$ws.Range("E27").Value = '  +1.59%  '
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").Value = '''7.36'
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("D31").Value = '''2.64'
$ws.Range("E31").Value = '  -3.13%  '
$ws.Range("D32").Value = '''2.13'
$ws.Range("E32").Value = '  -4.94%  '
$ws.Range("D33").Value = '''27.00'
$ws.Range("E33").Value = '  -4.99%  '
$ws.Range("E34").Value = '  -5.09%  '
$ws.Range("D35").Value = '''0.0₃0811'
$ws.Range("E35").Value = '  -5.84%  '
$ws.Range("E36").Value = '  -2.73%  '
$ws.Range("D37").Value = '''5.95'
$ws.Range("E37").Value = '  -4.27%  '
$ws.Range("D38").Value = '''3.22'
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("D39").Value = '''2.19'
$ws.Range("E39").Value = '  -5.75%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '''50.39'
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").Value = '''9.15'
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").Value = '''432.72'
$ws.Range("E42").Value = '  -7.48%  '
$ws.Range("D43").Value = '''0.283'
$ws.Range("E43").Value = '  -4.14%  '
$ws.Range("D44").Value = '''40.61'
$ws.Range("E44").Value = '  +1.90%  '
$ws.Range("D45").Value = '''0.111'
$ws.Range("E45").Value = '  +2.96%  '
$ws.Range("D46").Value = '''0.0357'
$ws.Range("E46").Value = '  -5.06%  '
$ws.Range("D47").Value = '''2.786.26'
$ws.Range("E47").Value = '  -4.26%  '
$ws.Range("D48").Value = '''129.38'
$ws.Range("E48").Value = '  -3.01%  '
$ws.Range("D50").Value = '''24.86'
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").Value = '''2.20'
$ws.Range("E51").Value = '  -2.42%  '
